# Model fitting based on data preprocessing
#
# Replace the raw traffic_volume strings (column D, rows 2-13) on the
# "Forecast Results" sheet with numeric 0, and refresh the computed
# MAE / MAPE / RMSE metrics on the "Metrics" sheet.

$wb = $excel.ActiveWorkbook

$wsForecast = $wb.Worksheets.Item("Forecast Results")
$wsMetrics  = $wb.Worksheets.Item("Metrics")

# Column D ("traffic_volume") holds inline-string blobs of raw data for
# rows 2 through 13; after preprocessing these collapse to a single
# numeric 0 value.
for ($row = 2; $row -le 13; $row++) {
    $wsForecast.Cells.Item($row, 4).Value = 0
}

# Updated model metrics resulting from the refit.
$wsMetrics.Range("A2").Value = 4.325770144240756
$wsMetrics.Range("B2").Value = 35.90012023593
$wsMetrics.Range("C2").Value = 5.56403091648521
